$wb = $excel.ActiveWorkbook

# Reverse the "*img" sheet names to "img*"
$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# The last sheet (formerly "eimg", now "imge") becomes the active sheet
$wb.Worksheets.Item("imge").Activate()
